$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, shifting existing rows 4-10 down to 5-11
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new data record
$ws.Range("A4").Value = 8
$ws.Range("B4").Value = "Terminal La Palmera de La Serena"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 44701
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 100112013
$ws.Range("G4").Value = "Alcachofa"
$ws.Range("H4").Value = "Española"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 19000
$ws.Range("L4").Value = 20000
$ws.Range("M4").Value = 19500
$ws.Range("N4").Value = "$/caja 30 unidades"
$ws.Range("O4").Value = "Provincia del Elquí"
$ws.Range("P4").Value = 650
$ws.Range("Q4").Value = 30
$ws.Range("R4").Value = "Hortaliza"
